# Auto-generated edit script
# Applies numeric corrections to market-price / profit columns (H-N)
# across multiple worksheets, per the target diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 299
$ws.Range("J2").Value = 299
$ws.Range("L2").Value = 299
$ws.Range("N2").Value = -525
$ws.Range("H9").Value = 152
$ws.Range("I9").Value = 159.75
$ws.Range("K9").Value = 159.75
$ws.Range("M9").Value = 9.25
$ws.Range("H70").Value = 852.5
$ws.Range("J70").Value = 852.5
$ws.Range("L70").Value = 2557.5
$ws.Range("N70").Value = -3097.5
$ws.Range("H73").Value = 852.5
$ws.Range("J73").Value = 852.5
$ws.Range("L73").Value = 2557.5
$ws.Range("N73").Value = -4429.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20230
$ws.Range("I32").Value = 17845
$ws.Range("K32").Value = 17845
$ws.Range("M32").Value = -17558
$ws.Range("H45").Value = 1096.6666
$ws.Range("I45").Value = 1096.6666
$ws.Range("K45").Value = 1096.6666
$ws.Range("M45").Value = -719.6666
$ws.Range("H61").Value = 3694
$ws.Range("J61").Value = 3694
$ws.Range("L61").Value = 3694
$ws.Range("N61").Value = -4118
$ws.Range("H136").Value = 3694
$ws.Range("J136").Value = 3694
$ws.Range("L136").Value = 11082
$ws.Range("N136").Value = -16182

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1192
$ws.Range("I94").Value = 1327.5
$ws.Range("J94").Value = 650
$ws.Range("K94").Value = 1327.5
$ws.Range("L94").Value = 650
$ws.Range("M94").Value = -876.5
$ws.Range("N94").Value = -1552

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 6200
$ws.Range("J4").Value = 8800
$ws.Range("L4").Value = 8800
$ws.Range("N4").Value = -9024
$ws.Range("H22").Value = 458.16666
$ws.Range("I22").Value = 437.5
$ws.Range("J22").Value = 499.5
$ws.Range("K22").Value = 437.5
$ws.Range("L22").Value = 499.5
$ws.Range("M22").Value = -87.5
$ws.Range("N22").Value = -1199.5
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").ClearContents()
$ws.Range("H42").Value = 10056
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("N49").ClearContents()
$ws.Range("H99").Value = 4028.7
$ws.Range("I99").Value = 3148
$ws.Range("K99").Value = 3148
$ws.Range("M99").Value = -1650
$ws.Range("H126").Value = 4028.7
$ws.Range("I126").Value = 3148
$ws.Range("K126").Value = 9444
$ws.Range("M126").Value = -6974
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 999.5
$ws.Range("J5").Value = 999
$ws.Range("L5").Value = 2997
$ws.Range("N5").Value = -3221
$ws.Range("H10").Value = 21.857143
$ws.Range("I10").Value = 13
$ws.Range("K10").Value = 39
$ws.Range("M10").Value = 100
$ws.Range("H26").Value = 26
$ws.Range("J26").Value = 101
$ws.Range("L26").Value = 303
$ws.Range("N26").Value = -879
$ws.Range("H34").Value = 694.44446
$ws.Range("I34").Value = 300
$ws.Range("J34").Value = 891.6667
$ws.Range("K34").Value = 900
$ws.Range("L34").Value = 2675.0001
$ws.Range("M34").Value = -816
$ws.Range("N34").Value = -2843.0001
$ws.Range("H39").Value = 1323.875
$ws.Range("J39").Value = 2247.75
$ws.Range("L39").Value = 6743.25
$ws.Range("N39").Value = -7331.25
$ws.Range("H55").Value = 200
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H104").Value = 3600
$ws.Range("J104").Value = 4000
$ws.Range("L104").Value = 12000
$ws.Range("N104").Value = -17242
$ws.Range("H131").Value = 2312.5
$ws.Range("J131").Value = 2346.6667
$ws.Range("L131").Value = 7040.000100000001
$ws.Range("N131").Value = -17120.0001
$ws.Range("H135").Value = 999.5
$ws.Range("J135").Value = 999
$ws.Range("L135").Value = 8991
$ws.Range("N135").Value = -14061
$ws.Range("H139").Value = 1777
$ws.Range("I139").Value = 1777
$ws.Range("K139").Value = 5331
$ws.Range("M139").Value = -191

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 3750
$ws.Range("I4").Value = 3750
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 3750
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -3638
$ws.Range("N4").ClearContents()
$ws.Range("H80").Value = 12325.272
$ws.Range("J80").Value = 17419.6
$ws.Range("L80").Value = 17419.6
$ws.Range("N80").Value = -19415.6
$ws.Range("H83").Value = 12325.272
$ws.Range("J83").Value = 17419.6
$ws.Range("L83").Value = 87098
$ws.Range("N83").Value = -97082
$ws.Range("H126").Value = 8012
$ws.Range("I126").Value = 8012
$ws.Range("K126").Value = 24036
$ws.Range("M126").Value = -21566

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1000
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H132").Value = 3490.889
$ws.Range("I132").Value = 3692.75
$ws.Range("K132").Value = 11078.25
$ws.Range("M132").Value = -8548.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 3483.25
$ws.Range("J4").Value = 4477
$ws.Range("L4").Value = 4477
$ws.Range("N4").Value = -4703
$ws.Range("H136").Value = 5108.222
$ws.Range("I136").Value = 3328.1667
$ws.Range("K136").Value = 9984.500100000001
$ws.Range("M136").Value = -7434.500100000001

Write-Host "Applied all cell updates."